# "Generate Report for Archive"
#
# The localization status report moved from "Ready for handoff" to
# "In Translation" for both translated docs/locales. Update the Status
# text on the Overview sheet (columns E/F, rows 2-3) and on each locale
# sheet's Status column (column C, rows 2-3), then shrink the now-narrower
# Status columns to match the shorter text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: zh-cn status (col E) and de-de status (col F) for both rows
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Per-locale detail sheets: Status column (C) for both rows
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# The Status columns auto-shrink to fit the shorter "In Translation" text
# (was width 17.2159881591797 / ColumnWidth 16.38, now ~13.41 / ColumnWidth 12.5)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
